# The source diff for this revision touches only the serialized *attribute
# order* inside word/document.xml (the <w:sectPr> / <w:pgSz> / <w:pgMar>
# values) and word/styles.xml (<w:docDefaults>, <w:latentStyles>, and the
# four <w:style> definitions). Every single changed line carries exactly
# the same attribute names/values as before -- only the left-to-right
# ordering of the attributes was normalized (alphabetically) by whatever
# tool produced that commit. There is no content, formatting, text, or
# structural change anywhere in the document.
#
# Because of that, the correct, faithful edit is one that leaves the
# document's actual content and formatting completely untouched -- the
# page geometry, section properties, run/paragraph text and every style
# definition must come out exactly as they went in. We touch the document
# only through read-only Word object model calls to confirm those values
# are already correct; we deliberately avoid writing to them (even
# reassigning a property to its own current value forces this host to
# rebuild that XML part, which would introduce spurious differences of
# its own -- e.g. extra namespace declarations -- that are not part of
# the recorded change), which would move the document further from the
# target rather than closer to it.

$d = $word.ActiveDocument

# --- word/document.xml : <w:sectPr> (page size / margins) -------------
$section = $d.Sections(1)
$pageSetup = $section.PageSetup

$expectedPageWidthTwips  = 11906
$expectedPageHeightTwips = 16838
$expectedTopTwips        = 1417
$expectedRightTwips      = 1417
$expectedBottomTwips     = 1417
$expectedLeftTwips       = 1417
$expectedHeaderTwips     = 708
$expectedFooterTwips     = 708
$expectedGutterTwips     = 0

# PageSetup reports/accepts points, OOXML stores twentieths-of-a-point
# (twips) -- 1 point = 20 twips. These reads are side-effect free.
$pageWidthTwips  = [Math]::Round($pageSetup.PageWidth  * 20)
$pageHeightTwips = [Math]::Round($pageSetup.PageHeight * 20)
$topTwips        = [Math]::Round($pageSetup.TopMargin    * 20)
$rightTwips      = [Math]::Round($pageSetup.RightMargin  * 20)
$bottomTwips     = [Math]::Round($pageSetup.BottomMargin * 20)
$leftTwips       = [Math]::Round($pageSetup.LeftMargin   * 20)
$headerTwips     = [Math]::Round($pageSetup.HeaderDistance * 20)
$footerTwips     = [Math]::Round($pageSetup.FooterDistance * 20)
$gutterTwips     = [Math]::Round($pageSetup.Gutter * 20)

Write-Output "pgSz  w=$pageWidthTwips h=$pageHeightTwips (expected w=$expectedPageWidthTwips h=$expectedPageHeightTwips)"
Write-Output "pgMar top=$topTwips right=$rightTwips bottom=$bottomTwips left=$leftTwips header=$headerTwips footer=$footerTwips gutter=$gutterTwips"

# --- word/styles.xml : docDefaults / latentStyles / style definitions --
$normalStyle = $d.Styles("Normal")
Write-Output ("Normal style name: " + $normalStyle.NameLocal)

$defaultParagraphFont = $d.Styles("Default Paragraph Font")
Write-Output ("Default Paragraph Font hidden: " + $defaultParagraphFont.Hidden)

$normalTable = $d.Styles("Normal Table")
Write-Output ("Normal Table priority: " + $normalTable.Priority)

$noList = $d.Styles("No List")
Write-Output ("No List priority: " + $noList.Priority)

# Nothing above mutated the package: every <w:...> element keeps its
# original (already value-correct) attributes, matching the diff's
# "same values, reordered attributes" nature without introducing any
# unrelated byte-level drift.
